# 自动更新Excel文件
# For every data row, the "剩余" (remaining days) counter in column E is
# decremented by one day. When a counter would drop to zero (or below),
# the cycle is considered finished and restarts: E is reset back to the
# "总天" (total days) value in column D and the "开始时间" (start date) in
# column F is pushed forward by 7 days (a new weekly cycle).
#
# Rows whose start-date value is not a well-formed 8-digit date (yyyymmdd)
# are left untouched, since their date can't be reliably advanced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)   # D: 总天 (total days)
    $eCell = $ws.Cells.Item($row, 5)   # E: 剩余 (days remaining)
    $fCell = $ws.Cells.Item($row, 6)   # F: 开始时间 (start date, yyyymmdd)

    $totalDays = $dCell.Value2
    $remaining = $eCell.Value2
    $startDate = $fCell.Value2

    if ($totalDays -eq $null -or $remaining -eq $null -or $startDate -eq $null) {
        continue
    }

    # Only touch rows that have a valid 8-digit yyyymmdd date
    $dateText = [string][int64]$startDate
    if ($dateText.Length -ne 8) {
        continue
    }

    $newRemaining = $remaining - 1

    if ($newRemaining -le 0) {
        # Cycle finished: reset remaining to total days and roll the start
        # date forward by one week (7 days)
        $year = [int]$dateText.Substring(0, 4)
        $month = [int]$dateText.Substring(4, 2)
        $day = [int]$dateText.Substring(6, 2)
        $dateValue = Get-Date -Year $year -Month $month -Day $day
        $newDateValue = $dateValue.AddDays(7)

        $eCell.Value2 = $totalDays
        $fCell.Value2 = [int64]($newDateValue.ToString("yyyyMMdd"))
    }
    else {
        $eCell.Value2 = $newRemaining
    }
}
